$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.264.95"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "2.595.88"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "578.07"
$ws.Range("E5").Value = "  +4.30%  "

$ws.Range("D6").Value = "142.88"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("D9").Value = "2.599.89"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").Value = "6.52"
$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").Value = "0.157"
$ws.Range("E12").Value = "  -1.83%  "

$ws.Range("D13").Value = "0.371"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").Value = "3.054.67"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "24.52"
$ws.Range("E15").Value = "  +6.61%  "

$ws.Range("D16").Value = "60.272.54"
$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").Value = "2.598.55"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "11.45"
$ws.Range("E19").Value = "  +10.01%  "

$ws.Range("D20").Value = "4.63"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").Value = "347.24"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +5.09%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").Value = "0.521"
$ws.Range("E24").Value = "  +7.84%  "

$ws.Range("D25").Value = "63.02"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").Value = "8.02"
$ws.Range("E28").Value = "  +7.28%  "

$ws.Range("D29").Value = "0.0₃0791"
$ws.Range("E29").Value = "  +3.33%  "

$ws.Range("E30").Value = "  +10.01%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  +4.73%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "164.82"
$ws.Range("E33").Value = "  +4.54%  "

$ws.Range("D34").Value = "19.43"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").Value = "4.29"
$ws.Range("E35").Value = "  +4.77%  "

$ws.Range("D36").Value = "0.983"
$ws.Range("E36").Value = "  +7.36%  "

$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  +8.10%  "

$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  +9.69%  "

$ws.Range("D39").Value = "38.07"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +6.52%  "

$ws.Range("D41").Value = "310.62"
$ws.Range("E41").Value = "  +7.31%  "

$ws.Range("D42").Value = "0.841"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").Value = "135.07"
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "0.0990"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("D46").Value = "5.02"
$ws.Range("E46").Value = "  +10.97%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "19.75"
$ws.Range("E47").Value = "  +4.82%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.603"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").Value = "0.0548"
$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("D50").Value = "20.05"
$ws.Range("E50").Value = "  +7.85%  "

$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  +2.73%  "
